# Updated solution for Tutorial 6
# Changes the date format from DD/MM/YYYY to DD-MM-YYYY in column A (rows 3-21)
# and updates the attendance counters (D, E, G, H) for a number of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New dash-separated dates for rows 3 through 21 (row -> date string)
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($row in $dates.Keys) {
    $cell = $ws.Range("A$row")
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$row]
}

# Rows where D (Total Attendance Count) -> 1, E (Real) -> 1, H (Absent) -> 0
$toPresent = @(4, 5, 6, 11, 12, 13, 14, 15, 16)
foreach ($row in $toPresent) {
    $ws.Range("D$row").Value = 1
    $ws.Range("E$row").Value = 1
    $ws.Range("H$row").Value = 0
}

# Row 20: D -> 1, G (Invalid) -> 1 (E stays 0, H stays 1)
$ws.Range("D20").Value = 1
$ws.Range("G20").Value = 1
